$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# read_through_names.py now also scrapes "height" and "weight" for each
# player. Insert two new columns before the existing "fantasy points"
# column (G), shifting fantasy points from G to I, then fill headers +
# the (constant, for this single player) height/weight values.
$ws.Columns("G").Insert()
$ws.Columns("G").Insert()

$ws.Range("G1").Value = "height"
$ws.Range("H1").Value = "weight"

$height = 6.166666666666667
$weight = 215

for ($row = 2; $row -le 8; $row++) {
    $ws.Cells.Item($row, 7).Value = $height
    $ws.Cells.Item($row, 8).Value = $weight
}
